# Add new weapons list
# - Items sheet: add 5 new weapon/item rows (Garrote, Poisoned Wind Globe, Scythe,
#   Spike Fist, Things Catcher) each hyperlinked to the "New Weapons" document.
# - Feats sheet: add 3 new weapon-related feat rows (Poisoned Wind Fumigator,
#   Spike Fist Brutalist, Things Wrangler) hyperlinked to the same document,
#   with the latter two sharing a merged hyperlink range.
# - Make the Feats sheet the active tab (it was MagicItems before).

$wb = $excel.ActiveWorkbook

$newWeaponsUrl = "https://editor.gmbinder.com/documents/edit/-N8HgfCb1-XMt8weFOvO"

# ---------------------------------------------------------------------------
# Items sheet (internal name "Items") -> 5 new rows
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items")

$itemRows = @(
    @{ Name = "Garrote";               Type = "Weapon" },
    @{ Name = "Poisoned Wind Globe";   Type = "Item" },
    @{ Name = "Scythe";                Type = "Weapon" },
    @{ Name = "Spike Fist";            Type = "Weapon" },
    @{ Name = "Things Catcher";        Type = "Weapon" }
)

$startRow = 6
for ($i = 0; $i -lt $itemRows.Count; $i++) {
    $r = $startRow + $i
    $row = $itemRows[$i]

    $items.Cells.Item($r, 1).Value = $row.Name
    $items.Cells.Item($r, 2).Value = $row.Type
    $items.Cells.Item($r, 3).Value = "New Weapons"
    $items.Hyperlinks.Add($items.Cells.Item($r, 3), $newWeaponsUrl) | Out-Null

    $items.Cells.Item($r, 4).Value = "Playtest Ready"
    $items.Cells.Item($r, 5).Value = "Not Released"
}

# Column C carries the "New Weapons" hyperlink text/style - copy formatting
# from an existing hyperlink cell so the new cells reuse the same style
# index instead of the ad-hoc one Hyperlinks.Add() would otherwise apply.
$items.Range("C2").Copy()
for ($i = 0; $i -lt $itemRows.Count; $i++) {
    $r = $startRow + $i
    $items.Cells.Item($r, 3).PasteSpecial(-4122)
}

# Extend the existing conditional formatting on column D to cover the new rows.
$itemsFcs = $items.Range("D1:D5").FormatConditions
$itemsNewRange = $items.Range("D1:D10")
for ($i = 1; $i -le $itemsFcs.Count; $i++) {
    $itemsFcs.Item($i).ModifyAppliesToRange($itemsNewRange)
}

# Column A grew wider to fit the new (longer) item names.
$items.Columns.Item(1).AutoFit()

$items.Range("D11").Select()

# ---------------------------------------------------------------------------
# Feats sheet -> 3 new rows
# ---------------------------------------------------------------------------
$feats = $wb.Worksheets.Item("Feats")

$featRows = @(
    @{ Name = "Poisoned Wind Fumigator" },
    @{ Name = "Spike Fist Brutalist" },
    @{ Name = "Things Wrangler" }
)

$featStartRow = 18
for ($i = 0; $i -lt $featRows.Count; $i++) {
    $r = $featStartRow + $i
    $row = $featRows[$i]

    $feats.Cells.Item($r, 1).Value = $row.Name
    $feats.Cells.Item($r, 2).Value = "Weapon"
    $feats.Cells.Item($r, 3).Value = "None"
    $feats.Cells.Item($r, 4).Value = "No"
    $feats.Cells.Item($r, 5).Value = "New Weapons"

    $feats.Cells.Item($r, 6).Value = "Playtest Ready"
    $feats.Cells.Item($r, 7).Value = "Not Released"
}

# Row 18 gets its own hyperlink; rows 19-20 share a single merged hyperlink
# with a custom display string, matching the source workbook's pattern for
# grouped "Independent Feats" style links.
$feats.Cells.Item(18, 5).Hyperlinks.Add($feats.Cells.Item(18, 5), $newWeaponsUrl) | Out-Null
$feats.Range("E19:E20").Hyperlinks.Add($feats.Range("E19:E20"), $newWeaponsUrl, "", "", "New Weapons") | Out-Null

# Column E carries the "New Weapons" hyperlink text/style - copy formatting
# from an existing hyperlink cell so the new cells reuse the same style
# index instead of the ad-hoc one Hyperlinks.Add() would otherwise apply.
$feats.Range("E2").Copy()
for ($i = 0; $i -lt $featRows.Count; $i++) {
    $r = $featStartRow + $i
    $feats.Cells.Item($r, 5).PasteSpecial(-4122)
}

# Column A grew wider to fit the new (longer) feat names.
$feats.Columns.Item(1).AutoFit()

$feats.Range("C23").Select()

# Feats becomes the active/visible tab (previously MagicItems was selected).
$feats.Activate()
